$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the filename parameter from "Kenya-2022" to "Namibia-2022"
$ws.Range("B8").Value = "Namibia-2022"

# Update the bounding-box coordinates for the new (Namibia) region
$ws.Range("B4").Value = 11    # Minimum longitude (deg)
$ws.Range("B5").Value = 26    # Maximum longitude (deg)
$ws.Range("B6").Value = -29   # Minimum latitude (deg)
$ws.Range("B7").Value = -16   # Maximum latitude (deg)

# Move the active selection to C5 (matches the saved view in the workbook)
$ws.Range("C5").Select()

$wb.Save()
